$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.129.95'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '2.480.00'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('D4').Value = "'1.00"
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = "'516.42"
$ws.Range('E5').Value = '  -0.60%  '
$ws.Range('D6').Value = "'131.45"
$ws.Range('E6').Value = '  +0.21%  '
$ws.Range('E7').Value = '  -0.37%  '
$ws.Range('D8').Value = "'0.553"
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').Value = '2.511.43'
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').Value = "'0.0974"
$ws.Range('E10').Value = '  -1.95%  '
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('D12').Value = "'5.22"
$ws.Range('E12').Value = '  -2.07%  '
$ws.Range('D13').Value = "'0.332"
$ws.Range('E13').Value = '  -3.37%  '
$ws.Range('D14').Value = '2.924.33'
$ws.Range('E14').Value = '  +0.24%  '
$ws.Range('D15').Value = '58.044.16'
$ws.Range('E15').Value = '  -0.04%  '
$ws.Range('D16').Value = "'22.16"
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = "'0.0000135"
$ws.Range('E17').Value = '  -1.11%  '
$ws.Range('D18').Value = '2.499.35'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = "'10.68"
$ws.Range('E19').Value = '  -1.51%  '
$ws.Range('D20').Value = "'320.62"
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = "'4.16"
$ws.Range('E21').Value = '  -0.41%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = "'5.99"
$ws.Range('E22').Value = '  +3.67%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  +0.13%  '
$ws.Range('D24').Value = "'63.15"
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = "'0.403"
$ws.Range('E25').Value = '  -1.73%  '
$ws.Range('B26').Value = 'Kaspa'
$ws.Range('C26').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D26').Value = "'0.161"
$ws.Range('E26').Value = '  +0.62%  '
$ws.Range('B27').Value = 'Binance-PegBSC-USD'
$ws.Range('C27').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D27').Value = "'0.990"
$ws.Range('E27').Value = '  -1.02%  '
$ws.Range('D28').Value = "'7.35"
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('D29').Value = "'169.66"
$ws.Range('E29').Value = '  +2.30%  '
$ws.Range('D30').Value = '0.0₃0744'
$ws.Range('E30').Value = '  -1.46%  '
$ws.Range('D31').Value = "'1.19"
$ws.Range('E31').Value = '  +1.62%  '
$ws.Range('B32').Value = 'Aptos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D32').Value = "'6.31"
$ws.Range('E32').Value = '  -0.04%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = "'1.69"
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D35').Value = "'0.994"
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').Value = "'18.02"
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('D37').Value = "'1.28"
$ws.Range('E37').Value = '  -2.84%  '
$ws.Range('D38').Value = "'3.95"
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('D39').Value = "'36.70"
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('E40').Value = '  -0.92%  '
$ws.Range('D41').Value = "'0.778"
$ws.Range('E41').Value = '  -1.51%  '
$ws.Range('D42').Value = "'277.25"
$ws.Range('E42').Value = '  +0.62%  '
$ws.Range('D43').Value = "'5.10"
$ws.Range('E43').Value = '  +1.60%  '
$ws.Range('D44').Value = "'3.41"
$ws.Range('E44').Value = '  -1.74%  '
$ws.Range('D45').Value = "'0.594"
$ws.Range('E45').Value = '  -0.03%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'122.90"
$ws.Range('E46').Value = '  -3.39%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = "'0.0917"
$ws.Range('E47').Value = '  +1.11%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'17.76"
$ws.Range('E48').Value = '  -1.23%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').Value = "'0.0494"
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').Value = "'0.0212"
$ws.Range('E50').Value = '  -0.80%  '
$ws.Range('D51').Value = "'16.87"
$ws.Range('E51').Value = '  -1.58%  '
